$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44846
$ws.Range("J2").Value = 7900
$ws.Range("D3").Value = 44204
$ws.Range("J3").Value = 7000
$ws.Range("D4").Value = 44875
$ws.Range("J4").Value = 7900
$ws.Range("D5").Value = 44883
$ws.Range("J5").Value = 9700
$ws.Range("D6").Value = 44187
$ws.Range("J6").Value = 12000
$ws.Range("D7").Value = 44214
$ws.Range("J7").Value = 7000
$ws.Range("D8").Value = 44186
$ws.Range("J8").Value = 10000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 30
$ws.Range("D9").Value = 44860
$ws.Range("J9").Value = 7900
$ws.Range("D10").Value = 44873
$ws.Range("D11").Value = 44902
$ws.Range("J11").Value = 7000
$ws.Range("D12").Value = 44162
$ws.Range("J12").Value = 7000
$ws.Range("D13").Value = 44876
$ws.Range("D14").Value = 44229
$ws.Range("J14").Value = 16000
$ws.Range("D15").Value = 44181
$ws.Range("J15").Value = 12000
$ws.Range("D16").Value = 44845
$ws.Range("J16").Value = 7900
$ws.Range("D17").Value = 44874
$ws.Range("J17").Value = 7900
$ws.Range("D18").Value = 44882
$ws.Range("D19").Value = 44230
$ws.Range("J19").Value = 16000
$ws.Range("D20").Value = 44231
$ws.Range("J20").Value = 12000
$ws.Range("D21").Value = 44602
$ws.Range("J21").Value = 12000
$ws.Range("D22").Value = 44602
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 2500
$ws.Range("O22").Value = "Provincia de Chacabuco"
$ws.Range("P22").Value = 25
$ws.Range("D23").Value = 44880
$ws.Range("J23").Value = 7900
$ws.Range("D24").Value = 44159
$ws.Range("J24").Value = 7000
$ws.Range("D25").Value = 44210
$ws.Range("J25").Value = 8800
$ws.Range("K25").Value = 2500
$ws.Range("M25").Value = 2750
$ws.Range("P25").Value = 28
$ws.Range("D26").Value = 44166
$ws.Range("J26").Value = 7000
$ws.Range("O26").Value = "Provincia de Chacabuco"
$ws.Range("D27").Value = 44167
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 7000
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = 3000
$ws.Range("O27").Value = "Provincia de Chacabuco"
$ws.Range("P27").Value = 30
$ws.Range("D28").Value = 44168
$ws.Range("D29").Value = 44901
$ws.Range("D30").Value = 44232
$ws.Range("J30").Value = 16000
$ws.Range("D31").Value = 44600
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 1300
$ws.Range("K31").Value = 3500
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = 3808
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 38
$ws.Range("D32").Value = 44209
$ws.Range("K32").Value = 2500
$ws.Range("M32").Value = 2750
$ws.Range("P32").Value = 28
$ws.Range("D33").Value = 44160
$ws.Range("D34").Value = 44245
$ws.Range("J34").Value = 9000
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("D35").Value = 44245
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 2500
$ws.Range("L35").Value = 2500
$ws.Range("M35").Value = 2500
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 25
$ws.Range("D36").Value = 44859
$ws.Range("J36").Value = 7900
$ws.Range("D37").Value = 44847
$ws.Range("D38").Value = 44189
$ws.Range("J38").Value = 16000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = 3000
$ws.Range("P38").Value = 30
$ws.Range("D39").Value = 44881
$ws.Range("J39").Value = 7900
$ws.Range("D40").Value = 44215
$ws.Range("J40").Value = 16000
$ws.Range("D42").Value = 44161
$ws.Range("J42").Value = 7000
$ws.Range("D43").Value = 44188
$ws.Range("J43").Value = 12000
